# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# Reordena los periodos de mora (columna E, filas 16-23) de orden
# descendente (2301,2212,...,2206) a orden ascendente (2206,...,2212,2301),
# y actualiza el Salario Basico (columna G) de 908526 a 1000000 para todas
# las filas. El valor de Valor Mora (columna F) de 33333 ahora corresponde
# al ultimo periodo (2301, fila 23) en vez del primero (fila 16).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$periods = @("2206", "2207", "2208", "2209", "2210", "2211", "2212", "2301")

for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = 16 + $i
    $ws.Cells.Item($row, 5).Value = $periods[$i]
    if ($row -eq 23) {
        $ws.Cells.Item($row, 6).Value = 33333
    } else {
        $ws.Cells.Item($row, 6).Value = 40000
    }
    $ws.Cells.Item($row, 7).Value = 1000000
}
